# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - refresh case counters for a handful of countries (Banglades, Afganistan,
#   Israel, Uzbekistan, Tailandia)
# - refresh the four low-count micro territories whose relative ranking
#   changed (Santa Sede / Islas Turcas y Caicos / Seychelles / Montserrat)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados a ..." footer -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 08:10"

# --- Banglades (row 21) -----------------------------------------------------
$ws.Range("D21").Value = 34027
$ws.Range("E21").Value = 55383

# --- Afganistan (row 42) -----------------------------------------------------
$ws.Range("B42").Value = 25633
$ws.Range("C42").Value = 106
$ws.Range("D42").Value = 5168
$ws.Range("E42").Value = 19981
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 484

# --- Israel (row 49) ---------------------------------------------------------
$ws.Range("B49").Value = 19338
$ws.Range("C49").Value = 101
$ws.Range("D49").Value = 15438
$ws.Range("E49").Value = 3598

# --- Uzbekistan (row 76) ------------------------------------------------------
$ws.Range("B76").Value = 5293
$ws.Range("C76").Value = 30
$ws.Range("E76").Value = 1255

# --- Tailandia (row 91) -------------------------------------------------------
$ws.Range("D91").Value = 2993
$ws.Range("E91").Value = 84

# --- Reordered micro territories (rows 208-211) ------------------------------
# Santa Sede now outranks Islas Turcas y Caicos, and Seychelles now outranks
# Montserrat, so the two pairs swap rows (country name + "Casos activos" +
# "Muertes").
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
